$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "30.963.11"
Set-TextValue $ws.Range("E2") "  +0.93%  "
Set-TextValue $ws.Range("D3") "1.952.99"
Set-TextValue $ws.Range("E3") "  -0.43%  "
Set-TextValue $ws.Range("D4") "1.002"
Set-TextValue $ws.Range("E4") "  +0.15%  "
Set-TextValue $ws.Range("D5") "244.79"
Set-TextValue $ws.Range("E5") "  -1.57%  "
Set-TextValue $ws.Range("E6") "  +0.03%  "
Set-TextValue $ws.Range("D7") "0.4881"
Set-TextValue $ws.Range("E7") "  +0.86%  "
Set-TextValue $ws.Range("D8") "0.2937"
Set-TextValue $ws.Range("E8") "  -0.62%  "
Set-TextValue $ws.Range("D9") "0.06801"
Set-TextValue $ws.Range("E9") "  +0.19%  "
Set-TextValue $ws.Range("D10") "19.16"
Set-TextValue $ws.Range("E10") "  -0.92%  "
Set-TextValue $ws.Range("D11") "106.90"
Set-TextValue $ws.Range("E11") "  -3.35%  "
Set-TextValue $ws.Range("D12") "1.953.28"
Set-TextValue $ws.Range("D13") "0.07788"
Set-TextValue $ws.Range("E13") "  +0.57%  "
Set-TextValue $ws.Range("D14") "5.422"
Set-TextValue $ws.Range("D15") "0.6982"
Set-TextValue $ws.Range("E15") "  +1.10%  "
Set-TextValue $ws.Range("D16") "279.10"
Set-TextValue $ws.Range("E16") "  -4.82%  "
Set-TextValue $ws.Range("D17") "30.983.96"
Set-TextValue $ws.Range("E17") "  +0.99%  "
Set-TextValue $ws.Range("D18") "13.16"
Set-TextValue $ws.Range("E18") "  -1.27%  "
Set-TextValue $ws.Range("D19") "0.000007658"
Set-TextValue $ws.Range("E19") "  -0.42%  "
Set-TextValue $ws.Range("D20") "2.205.94"
Set-TextValue $ws.Range("E20") "  -0.40%  "
Set-TextValue $ws.Range("E21") "  +0.05%  "
Set-TextValue $ws.Range("D22") "5.468"
Set-TextValue $ws.Range("E22") "  -3.34%  "
Set-TextValue $ws.Range("D23") "1.002"
Set-TextValue $ws.Range("E23") "  +0.12%  "
Set-TextValue $ws.Range("D24") "6.466"
Set-TextValue $ws.Range("E24") "  -2.26%  "
Set-TextValue $ws.Range("D25") "9.800"
Set-TextValue $ws.Range("E25") "  -0.95%  "
Set-TextValue $ws.Range("D26") "169.51"
Set-TextValue $ws.Range("E26") "  -0.52%  "
Set-TextValue $ws.Range("D27") "19.93"
Set-TextValue $ws.Range("E27") "  -1.05%  "
Set-TextValue $ws.Range("D28") "2.190"
Set-TextValue $ws.Range("E28") "  -0.47%  "
Set-TextValue $ws.Range("E29") "  -1.62%  "
Set-TextValue $ws.Range("D30") "1.410"
Set-TextValue $ws.Range("E30") "  -2.05%  "
Set-TextValue $ws.Range("D31") "1.576"
Set-TextValue $ws.Range("E31") "  -1.52%  "
Set-TextValue $ws.Range("D32") "4.606"
Set-TextValue $ws.Range("E32") "  -1.59%  "
Set-TextValue $ws.Range("D33") "4.430"
Set-TextValue $ws.Range("E33") "  -0.70%  "
Set-TextValue $ws.Range("D34") "0.04927"
Set-TextValue $ws.Range("E34") "  -3.70%  "
Set-TextValue $ws.Range("D35") "0.7618"
Set-TextValue $ws.Range("E35") "  -2.30%  "
Set-TextValue $ws.Range("D36") "1.165"
Set-TextValue $ws.Range("E36") "  -0.89%  "
Set-TextValue $ws.Range("D37") "2.731"
Set-TextValue $ws.Range("E37") "  -0.12%  "
Set-TextValue $ws.Range("D38") "0.02000"
Set-TextValue $ws.Range("E38") "  -2.89%  "
Set-TextValue $ws.Range("D39") "2.705"
Set-TextValue $ws.Range("E39") "  -0.21%  "
Set-TextValue $ws.Range("D40") "6.477"
Set-TextValue $ws.Range("E40") "  +5.70%  "
Set-TextValue $ws.Range("D41") "2.116"
Set-TextValue $ws.Range("E41") "  +2.06%  "
Set-TextValue $ws.Range("D42") "73.79"
Set-TextValue $ws.Range("E42") "  +4.79%  "
Set-TextValue $ws.Range("D43") "0.8845"
Set-TextValue $ws.Range("E43") "  +0.99%  "
Set-TextValue $ws.Range("D46") "8.089"
Set-TextValue $ws.Range("E46") "  +9.35%  "
Set-TextValue $ws.Range("E47") "  -0.05%  "
Set-TextValue $ws.Range("D48") "989.53"
Set-TextValue $ws.Range("E48") "  +8.93%  "
Set-TextValue $ws.Range("D49") "0.1257"
Set-TextValue $ws.Range("E49") "  -2.26%  "
Set-TextValue $ws.Range("D50") "9.250"
Set-TextValue $ws.Range("E50") "  -1.02%  "
Set-TextValue $ws.Range("D51") "0.2564"
Set-TextValue $ws.Range("E51") "  +1.94%  "

# Row 44/45: Quant and TheSandbox swap positions with updated values
Set-TextValue $ws.Range("B44") "Quant"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D44") "109.07"
Set-TextValue $ws.Range("E44") "  -1.89%  "

Set-TextValue $ws.Range("B45") "TheSandbox"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D45") "0.4442"
Set-TextValue $ws.Range("E45") "  -0.43%  "
